$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two rows were removed from the source data (RM 232, then SC 92), which
# shifts everything below them up by one row each time. After the first
# delete, SC 92 (originally row 28) has shifted up to row 27.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Remaining per-cell value edits (newly-imputed values filled in, and a
# few values newly blanked out), against the now-shifted row numbers.
$ws.Range("E2").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("D6").Value = -14.2
$ws.Range("F6").Value = 16.43
$ws.Range("D8").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("F14").Value = 17.76
$ws.Range("D18").Value = -15.2
$ws.Range("D20").ClearContents()
$ws.Range("F21").Value = 16.58
$ws.Range("F22").Value = 16.81
$ws.Range("D23").Value = -13.9
$ws.Range("D25").ClearContents()
$ws.Range("F26").ClearContents()
$ws.Range("B27").Value = -20.4
$ws.Range("F27").ClearContents()
$ws.Range("B28").ClearContents()
$ws.Range("F28").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("B30").Value = -19.7
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F31").Value = 17.18
$ws.Range("B32").ClearContents()
